$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.701.66"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "1.632.78"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "218.01"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("E6").Value = "  -1.57%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -1.39%  "
$ws.Range("D11").Value = "0.0843"
$ws.Range("D12").Value = "1.860.97"
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D13").Value = "1.626.71"
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("D15").Value = "0.522"
$ws.Range("E15").Value = "  -2.11%  "
$ws.Range("E16").Value = "  -2.59%  "
$ws.Range("D17").Value = "26.687.10"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "0.0₃0723"
$ws.Range("E18").Value = "  -2.84%  "
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").Value = "210.98"
$ws.Range("E20").Value = "  -3.06%  "
$ws.Range("D21").Value = "4.31"
$ws.Range("E21").Value = "  -1.63%  "
$ws.Range("E22").Value = "  -2.22%  "
$ws.Range("D23").Value = "2.33"
$ws.Range("E23").Value = "  -3.50%  "
$ws.Range("D24").Value = "9.17"
$ws.Range("E24").Value = "  -3.14%  "
$ws.Range("D25").Value = "146.86"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("E27").Value = "  -2.64%  "
$ws.Range("E28").Value = "  -3.26%  "
$ws.Range("D29").Value = "15.54"
$ws.Range("E29").Value = "  -1.86%  "
$ws.Range("D30").Value = "0.0501"
$ws.Range("E30").Value = "  -4.31%  "
$ws.Range("E31").Value = "  +0.61%  "
$ws.Range("D32").Value = "3.37"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").Value = "2.95"
$ws.Range("E33").Value = "  -2.58%  "
$ws.Range("D34").Value = "1.262.17"
$ws.Range("E34").Value = "  -1.23%  "
$ws.Range("E35").Value = "  -2.08%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("E37").Value = "  -3.36%  "
$ws.Range("D38").Value = "0.523"
$ws.Range("E38").Value = "  -3.50%  "
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("E40").Value = "  -3.95%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "2.15"
$ws.Range("E42").Value = "  -4.46%  "
$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D43").Value = "1.771.66"
$ws.Range("E43").Value = "  -1.53%  "
$ws.Range("E44").Value = "  -3.24%  "
$ws.Range("D45").Value = "91.28"
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("D46").Value = "59.62"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").Value = "1.56"
$ws.Range("E47").Value = "  -3.78%  "
$ws.Range("D48").Value = "0.0516"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").Value = "0.407"
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("D51").Value = "0.0955"
$ws.Range("E51").Value = "  -2.73%  "
